$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TestScenario_1 block): Approved/Rejected column flips from "Approved" to
# "Rejected", and a reason is supplied in the previously-empty ReasonToReject column.
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "test"

# Row 19 (TestScenario_2 block): same Approved -> Rejected flip (shared string reused),
# plus a reason in the ReasonToReject column.
$ws.Range("I19").Value = "Rejected"
$ws.Range("J19").Value = "eerere"

# Update the active selection to match where the editor left off.
$null = $ws.Range("I16").Select()
